# This script rearranges the data rows 2-11 of the active sheet.
# Columns D, L, M, N, O, P, Q, R, S, T are permuted across rows (the
# other columns - A,B,C,E,F,G,H,I,J,K - are identical for every row
# already, so they do not need to change).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of destination row -> source row (values to copy from the
# "before" layout into the destination row).
$mapping = @{
    2  = 6
    3  = 10
    4  = 7
    5  = 2
    6  = 8
    7  = 3
    8  = 11
    9  = 5
    10 = 9
    11 = 4
}

# Columns whose values move as part of the permutation.
$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Snapshot the current ("before") values for the columns involved, for
# every row, before we start overwriting anything.
$snapshot = @{}
foreach ($row in 2..11) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowData
}

# Now write the permuted values into each destination row.
foreach ($destRow in 2..11) {
    $srcRow = $mapping[$destRow]
    $srcData = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $srcData[$col]
    }
}
